# "Submitted timeline week 8"
# Fills in the Week 8 timesheet with the week's entries (Feb 24 - Mar 1, 2018)
# and makes "Week 8" the active/selected sheet (previously "Week 7" was active).
#
# Dates/times are written as raw serial numbers (this workbook uses the
# date1904 base) with an explicit NumberFormat, rather than locale date/time
# strings, so the stored values line up with the 1904 epoch the workbook
# already uses for every other week's entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 8")

# --- Row 2: Sat 2/24, 10:00 PM - 12:00 AM, 2 hrs ---
$ws.Range("A2").Value = 41693
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Value = 0.91666666666666663
$ws.Range("B2").NumberFormat = "h:mm AM/PM"
$ws.Range("C2").Value = 0
$ws.Range("C2").NumberFormat = "h:mm AM/PM"
$ws.Range("D2").Value = "Created dropdown sign in on navbar, adjusted CSS"
$ws.Range("E2").Value = 2
$ws.Rows.Item(2).RowHeight = 18

# --- Row 3: Sun 2/25, 10:00 PM - 12:00 AM, 2 hrs ---
$ws.Range("A3").Value = 41694
$ws.Range("A3").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Value = 0.91666666666666663
$ws.Range("B3").NumberFormat = "h:mm AM/PM"
$ws.Range("C3").Value = 0
$ws.Range("C3").NumberFormat = "h:mm AM/PM"
$ws.Range("D3").Value = "Created sign up page, linked page in navigation"
$ws.Range("E3").Value = 2
$ws.Rows.Item(3).RowHeight = 18

# --- Row 4: Mon 2/26, 9:00 AM - 10:00 AM, 1 hr ---
$ws.Range("A4").Value = 41695
$ws.Range("A4").NumberFormat = "mm-dd-yy"
$ws.Range("B4").Value = 0.375
$ws.Range("B4").NumberFormat = "h:mm AM/PM"
$ws.Range("C4").Value = 0.41666666666666669
$ws.Range("C4").NumberFormat = "h:mm AM/PM"
$ws.Range("D4").Value = "Created and tested form to add new user which redirects to index, added roles and policies to users, created ""manage users"" page for admin role"
$ws.Range("D4").WrapText = $true
$ws.Range("E4").Value = 1
$ws.Rows.Item(4).RowHeight = 39

# --- Row 5: Mon 2/26, 4:30 PM - 5:30 PM, 1 hr ---
$ws.Range("A5").Value = 41695
$ws.Range("A5").NumberFormat = "mm-dd-yy"
$ws.Range("B5").Value = 0.6875
$ws.Range("B5").NumberFormat = "h:mm AM/PM"
$ws.Range("C5").Value = 0.72916666666666663
$ws.Range("C5").NumberFormat = "h:mm AM/PM"
$ws.Range("D5").Value = "Created sign in views and functions in UserController"
$ws.Range("E5").Value = 1
$ws.Rows.Item(5).RowHeight = 18

# --- Row 6: Tue 2/27, 10:00 AM - 12:00 PM, 2 hrs ---
$ws.Range("A6").Value = 41696
$ws.Range("A6").NumberFormat = "mm-dd-yy"
$ws.Range("B6").Value = 0.41666666666666669
$ws.Range("B6").NumberFormat = "h:mm AM/PM"
$ws.Range("C6").Value = 0.5
$ws.Range("C6").NumberFormat = "h:mm AM/PM"
$ws.Range("D6").Value = "Sign up auto login, redirects to user account, authorization on guest (cannot view account unless logged in), created cart model to add products to cart"
$ws.Range("D6").WrapText = $true
$ws.Range("E6").Value = 2
$ws.Rows.Item(6).RowHeight = 39

# --- Row 7: Tue 2/27, 4:00 PM - Wed 5:30 AM (overnight), 1.5 hrs ---
$ws.Range("A7").Value = 41696
$ws.Range("A7").NumberFormat = "mm-dd-yy"
$ws.Range("B7").Value = 0.66666666666666663
$ws.Range("B7").NumberFormat = "h:mm AM/PM"
$ws.Range("C7").Value = 0.22916666666666666
$ws.Range("C7").NumberFormat = "h:mm"
$ws.Range("D7").Value = "Added images to products table, successful display of images with products on products page"
$ws.Range("D7").WrapText = $true
$ws.Range("E7").Value = 1.5
$ws.Rows.Item(7).RowHeight = 26

# --- Row 8: Thu 3/1, 10:00 AM - 12:00 PM, 2 hrs ---
$ws.Range("A8").Value = 41698
$ws.Range("A8").NumberFormat = "mm-dd-yy"
$ws.Range("B8").Value = 0.41666666666666669
$ws.Range("B8").NumberFormat = "h:mm AM/PM"
$ws.Range("C8").Value = 0.5
$ws.Range("C8").NumberFormat = "h:mm AM/PM"
$ws.Range("D8").Value = "Successful add to cart, view cart, and started on checkout"
$ws.Range("D8").WrapText = $true
$ws.Range("E8").Value = 2
$ws.Rows.Item(8).RowHeight = 18

# Make "Week 8" the active sheet with A9 selected (next blank entry row),
# matching the workbook's saved view state after submitting this week.
$ws.Activate()
$ws.Range("A9").Select()
